$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the existing row 284, shifting rows
# 284-402 down to 285-403 (matches the diff: old row N becomes new row N+1
# for N in 284..402, and a fresh row 284 is introduced).
$ws.Range("A284").EntireRow.Insert()

# Populate the newly inserted row 284 with the new Perejil price entry.
$ws.Range("A284").Value = 9
$ws.Range("B284").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C284").Value = 'Metropolitana'
$ws.Range("D284").Value = 44755
$ws.Range("E284").Value = 13
$ws.Range("F284").Value = 100112044
$ws.Range("G284").Value = 'Perejil'
$ws.Range("H284").Value = 'Sin especificar'
$ws.Range("I284").Value = 'Primera'
$ws.Range("J284").Value = 52
$ws.Range("K284").Value = 19000
$ws.Range("L284").Value = 20000
$ws.Range("M284").Value = 19500
$ws.Range("N284").Value = '$/docena de atados'
$ws.Range("O284").Value = 'Región Metropolitana'
$ws.Range("P284").Value = 6500
$ws.Range("Q284").Value = 3
$ws.Range("R284").Value = 'Hortaliza'
